$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cmaes")

# Rename knob label (row 7, column A): knobs.knob_ids -> knobs.knob_link_ids
$ws.Range("A7").Value = "knobs.knob_link_ids"

# Update beats_parameters.SIM_DT (row 16) for columns B, C, D (E stays [5])
$ws.Range("B16").Value = "[4]"
$ws.Range("C16").Value = "[4]"
$ws.Range("D16").Value = "[4]"

# Update beats_parameters.DURATION (row 14, col B) and beats_parameters.OUTPUT_DT (row 15, col B)
$ws.Range("B14").Value = "[86400]"
$ws.Range("B15").Value = "[300]"

# Setting .Value above resets the cells' "quote prefix" number format (style index 1)
# back to the default style. Restore it by copying the format from an untouched sibling
# cell in the same row that still carries that style.
$ws.Range("E14").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

$ws.Range("E15").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$ws.Range("E16").Copy() | Out-Null
$ws.Range("B16:D16").PasteSpecial(-4122) | Out-Null

# Move the active cell selection from B7 to B6
$ws.Range("B6").Select()
